$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "Done"
$ws.Range("C6").Value = "Done"
$ws.Range("D6").Value = "Done"
$ws.Range("E6").Value = "ETL"

$ws.Range("B7").Value = "Done"
$ws.Range("C7").Value = "Done"
$ws.Range("D7").Value = "Done"
$ws.Range("E7").Value = "Data Cleaning"

$ws.Range("E7").Select()
